$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1999170.6
$ws.Range("I33").Value = 2713067
$ws.Range("K33").Value = 2713067
$ws.Range("M33").Value = -2712838

$ws.Range("H41").Value = 724.6667
$ws.Range("J41").Value = 938.625
$ws.Range("L41").Value = 938.625
$ws.Range("N41").Value = -1818.625

$ws.Range("H64").Value = 3554
$ws.Range("I64").Value = 3541.6667
$ws.Range("K64").Value = 3541.6667
$ws.Range("M64").Value = -3293.6667

$ws.Range("H67").Value = 3554
$ws.Range("I67").Value = 3541.6667
$ws.Range("K67").Value = 3541.6667
$ws.Range("M67").Value = -2683.6667

$ws.Range("H86").Value = 2029.4166
$ws.Range("I86").Value = 2225.158
$ws.Range("J86").Value = 1285.6
$ws.Range("K86").Value = 2225.158
$ws.Range("L86").Value = 1285.6
$ws.Range("M86").Value = -1102.158
$ws.Range("N86").Value = -3531.6

$ws.Range("H89").Value = 2029.4166
$ws.Range("I89").Value = 2225.158
$ws.Range("J89").Value = 1285.6
$ws.Range("K89").Value = 11125.79
$ws.Range("L89").Value = 6428
$ws.Range("M89").Value = -5509.789999999999
$ws.Range("N89").Value = -17660

$ws.Range("H96").Value = 1210.3334
$ws.Range("I96").Value = 554.3333
$ws.Range("K96").Value = 1662.9999
$ws.Range("M96").Value = -289.9999

$ws.Range("H98").Value = 1404.8572
$ws.Range("I98").Value = 1358.4445
$ws.Range("K98").Value = 1358.4445
$ws.Range("M98").Value = 139.5554999999999

$ws.Range("H106").Value = 20355.562
$ws.Range("I106").Value = 16999.154
$ws.Range("J106").Value = 34900
$ws.Range("K106").Value = 16999.154
$ws.Range("L106").Value = 34900
$ws.Range("M106").Value = -16368.154
$ws.Range("N106").Value = -36162

$ws.Range("H122").Value = 1404.8572
$ws.Range("I122").Value = 1358.4445
$ws.Range("K122").Value = 4075.3335
$ws.Range("M122").Value = -1625.3335

$ws.Range("H132").Value = 6900425.5
$ws.Range("I132").Value = 9093856
$ws.Range("K132").Value = 27281568
$ws.Range("M132").Value = -27279038

$ws.Range("H137").Value = 8897.405000000001
$ws.Range("I137").Value = 6970.5
$ws.Range("K137").Value = 20911.5
$ws.Range("M137").Value = -18361.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H82").Value = 37499.5
$ws.Range("J82").Value = 37499.5
$ws.Range("L82").Value = 37499.5
$ws.Range("N82").Value = -38221.5

$ws.Range("H85").Value = 37499.5
$ws.Range("J85").Value = 37499.5
$ws.Range("L85").Value = 37499.5
$ws.Range("N85").Value = -39995.5

$ws.Range("H94").Value = 49497.5
$ws.Range("J94").Value = 49497.5
$ws.Range("L94").Value = 49497.5
$ws.Range("N94").Value = -51299.5

$ws.Range("H132").Value = 7223.9575
$ws.Range("I132").Value = 5827.582
$ws.Range("K132").Value = 17482.746
$ws.Range("M132").Value = -14952.746

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H99").Value = 3375
$ws.Range("I99").Value = 3375
$ws.Range("K99").Value = 3375
$ws.Range("M99").Value = -1877

$ws.Range("H107").Value = 5500
$ws.Range("I107").Value = 3689.3
$ws.Range("J107").Value = 8517.833000000001
$ws.Range("K107").Value = 3689.3
$ws.Range("L107").Value = 8517.833000000001
$ws.Range("M107").Value = -1769.3
$ws.Range("N107").Value = -12357.833

$ws.Range("H134").Value = 6736.1914
$ws.Range("I134").Value = 3968.4443
$ws.Range("K134").Value = 11905.3329
$ws.Range("M134").Value = -9370.332900000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2050048.9
$ws.Range("I31").Value = 2841560
$ws.Range("J31").Value = 1431.4117
$ws.Range("K31").Value = 2841560
$ws.Range("L31").Value = 1431.4117
$ws.Range("M31").Value = -2841265
$ws.Range("N31").Value = -2021.4117

$ws.Range("H34").Value = 2050048.9
$ws.Range("I34").Value = 2841560
$ws.Range("J34").Value = 1431.4117
$ws.Range("K34").Value = 2841560
$ws.Range("L34").Value = 1431.4117
$ws.Range("M34").Value = -2841358
$ws.Range("N34").Value = -1835.4117

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H132").Value = 27026.684
$ws.Range("I132").Value = 21803.45
$ws.Range("J132").Value = 32830.277
$ws.Range("K132").Value = 65410.35000000001
$ws.Range("L132").Value = 98490.83100000001
$ws.Range("M132").Value = -62880.35000000001
$ws.Range("N132").Value = -103550.831

$ws.Range("H134").Value = 11838.323
$ws.Range("I134").Value = 9277.4
$ws.Range("J134").Value = 15496.786
$ws.Range("K134").Value = 27832.2
$ws.Range("L134").Value = 46490.358
$ws.Range("M134").Value = -25297.2
$ws.Range("N134").Value = -51560.358

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2799.8
$ws.Range("I11").Value = 2799.8
$ws.Range("K11").Value = 8399.400000000001
$ws.Range("M11").Value = -8259.400000000001

$ws.Range("H17").Value = 1195.9231
$ws.Range("I17").Value = 410
$ws.Range("J17").Value = 1687.125
$ws.Range("K17").Value = 1230
$ws.Range("L17").Value = 5061.375
$ws.Range("M17").Value = -1061
$ws.Range("N17").Value = -5399.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H132").Value = 9200.419
$ws.Range("I132").Value = 5257.914
$ws.Range("K132").Value = 15773.742
$ws.Range("M132").Value = -13243.742

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H81").Value = 59999
$ws.Range("J81").Value = 59999
$ws.Range("L81").Value = 59999
$ws.Range("N81").Value = -61995

$ws.Range("H84").Value = 59999
$ws.Range("J84").Value = 59999
$ws.Range("L84").Value = 179997
$ws.Range("N84").Value = -189981

$ws.Range("H98").Value = 50000
$ws.Range("J98").Value = 50000
$ws.Range("L98").Value = 50000
$ws.Range("N98").Value = -55990

$ws.Range("H122").Value = 4202.8184
$ws.Range("I122").Value = 3391.4
$ws.Range("K122").Value = 10174.2
$ws.Range("M122").Value = -7724.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 25640.824
$ws.Range("J54").Value = 25640.824
$ws.Range("L54").Value = 25640.824
$ws.Range("N54").Value = -26680.824

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H63").Value = 39999
$ws.Range("J63").Value = 39999
$ws.Range("L63").Value = 39999
$ws.Range("N63").Value = -41247

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H66").Value = 39999
$ws.Range("J66").Value = 39999
$ws.Range("L66").Value = 119997
$ws.Range("N66").Value = -126237

$ws.Range("H74").Value = 18999
$ws.Range("I74").Value = 18999
$ws.Range("K74").Value = 18999
$ws.Range("M74").Value = -18063

$ws.Range("H77").Value = 18999
$ws.Range("I77").Value = 18999
$ws.Range("K77").Value = 56997
$ws.Range("M77").Value = -52317

$ws.Range("H92").Value = 47500
$ws.Range("J92").Value = 47500
$ws.Range("L92").Value = 47500
$ws.Range("N92").Value = -52492

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H98").Value = 63959
$ws.Range("J98").Value = 63959
$ws.Range("L98").Value = 63959
$ws.Range("N98").Value = -69949

$ws.Range("H136").Value = 2421.6
$ws.Range("I136").Value = 2155.6743
$ws.Range("K136").Value = 6467.0229
$ws.Range("M136").Value = -3917.0229
